# Apply "CDS Study filter fixes" commit:
#  - Replace the Participants-tab query (B2) with the corrected version that
#    matches participant->study (instead of study<-participant), adds the
#    genomic_info hop, and sorts the collected sample ids.
#  - Replace the Files-tab query (B4) with the corrected version that matches
#    file->study (instead of study<-participant<-sample<-file), adds the
#    genomic_info hop, and renames the "Subject ID" column to "Participant ID".
#  - Row heights on rows 2 and 4 grow to fit the new (taller) query text.
#  - Selection moves to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Human Tumor Atlas Network (HTAN) primary sequencing data"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

$newFileQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Human Tumor Atlas Network (HTAN) primary sequencing data"]
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name
Limit 100
'@

# ParticipantsTab row (row 2): update the query text in column B.
$ws.Range("B2").Value = $newParticipantQuery

# FilesTab row (row 4): update the query text in column B.
$ws.Range("B4").Value = $newFileQuery

# Row heights grow because the replacement query text wraps across more
# lines than the text it replaced.
$ws.Rows.Item(2).RowHeight = 283.5
$ws.Rows.Item(4).RowHeight = 267.75

# Selection ends on E4 (matches the saved sheetView selection in the diff).
$ws.Range("E4").Select() | Out-Null
